$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # top
$c1.Borders.Item(9).LineStyle = 1   # bottom
Write-Host "C1 done"

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(10).LineStyle = 1  # right
$d1.Borders.Item(8).LineStyle = 1   # top
$d1.Borders.Item(9).LineStyle = 1   # bottom
Write-Host "D1 done"
